$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60; this shifts the existing rows 60..83 down to 61..84,
# which already reproduces the "rotation" seen across the rest of the diff (each old
# row's data reappears one row lower, keeping its original contents).
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new weekly record.
$ws.Cells.Item(60, 1).Value2 = 11
$ws.Cells.Item(60, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(60, 3).Value2 = "Bíobío"
$ws.Cells.Item(60, 4).Value2 = 45135
$ws.Cells.Item(60, 5).Value2 = 8
$ws.Cells.Item(60, 6).Value2 = 100112043
$ws.Cells.Item(60, 7).Value2 = "Pepino dulce"
$ws.Cells.Item(60, 8).Value2 = "Sin especificar"
$ws.Cells.Item(60, 9).Value2 = "Primera"
$ws.Cells.Item(60, 10).Value2 = 100
$ws.Cells.Item(60, 11).Value2 = 17000
$ws.Cells.Item(60, 12).Value2 = 18000
$ws.Cells.Item(60, 13).Value2 = 17500
$ws.Cells.Item(60, 14).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(60, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(60, 16).Value2 = 972
$ws.Cells.Item(60, 17).Value2 = 18
$ws.Cells.Item(60, 18).Value2 = "Hortaliza"
